$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9 - shifts existing rows 9..29 down to 10..30
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new data record
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Terminal La Palmera de La Serena"
$ws.Range("C9").Value = "Coquimbo"
$ws.Range("D9").Value = 44659
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100101
$ws.Range("H9").Value = "Berries"
$ws.Range("I9").Value = 100101001
$ws.Range("J9").Value = "Arándano (blue)"
$ws.Range("K9").Value = "Sin especificar"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 200
$ws.Range("N9").Value = 6500
$ws.Range("O9").Value = 7000
$ws.Range("P9").Value = 6750
$ws.Range("Q9").Value = '$/bandeja 2 kilos'
$ws.Range("R9").Value = "Provincia de Linares"
$ws.Range("S9").Value = 3375
$ws.Range("T9").Value = 2
